# Project Logbook.xlsx - add new logbook entry row (row 40):
#   Start Date / Finish Date : 12/11/2017 (serial 43080)
#   Member                   : Antonio Vazquez
#   Activity                 : CAN code issues with the initialization solved
#                               (still present some problems with the timer)
# Also move the active selection to the new row (D40), matching the saved
# view state captured in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New logbook row --------------------------------------------------
$ws.Range("A40").Value = 43080
$ws.Range("A40").NumberFormat = "mm-dd-yy"

# Share the exact same cell style between A40 and B40 (copy/paste the
# format instead of re-applying NumberFormat so only one new style entry
# is minted instead of two identical ones).
$ws.Range("A40").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$ws.Range("B40").Value = 43080

$ws.Range("C40").Value = "Antonio Vazquez"
$ws.Range("D40").Value = "CAN code issues with the initialization solved (still present some problems with the timer)"

# --- Selection / view ---------------------------------------------------
$ws.Range("D40").Select()
